# Facebook ad processing data cleanup:
# - Replace "FIX THIS" placeholders with real Facebook page_id values
# - Clear malformed / superseded page_id_2 values
# - Move a page_id that had landed in the wrong column back to page_id (B)
# - Update the saved selection / scroll position to reflect where the
#   author was last working in the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# bruce rauner (row 50) - page_id was "FIX THIS"
$ws.Range("B50").Value = "213918568751224"

# dean heller (row 97) - page_id_2 was a malformed placeholder; clear it
$ws.Range("C97").Value = $null

# deb fischer (row 98) - the correct page_id belongs in column B;
# the malformed page_id_2 in column C is cleared
$ws.Range("B98").Value = "109592402468562"
$ws.Range("C98").Value = $null

# deidre dejear (row 100) - fill in the newly-found page_id
$ws.Range("B100").Value = "1490813904322906"

# karl dean (row 213) - page_id was "FIX THIS"
$ws.Range("B213").Value = "191055010919565"

# kay ivey (row 220) - page_id_2 was "FIX THIS"; clear it
$ws.Range("C220").Value = $null

# kevin o'connor (row 224) - page_id was "FIX THIS"
$ws.Range("B224").Value = "104346771140701"

# kyrsten sinema (row 231) - page_id was "FIX THIS"
$ws.Range("B231").Value = "52563647525"

# mary throne (row 269) - page_id was "FIX THIS"
$ws.Range("B269").Value = "1409012545881430"

# nicole galloway (row 298) - page_id was "FIX THIS"
$ws.Range("B298").Value = "149366465128820"

# shelley lenz (row 344) - page_id was "FIX THIS"
$ws.Range("B344").Value = "100110638215668"

# Reflect the author's last scroll/selection position in the sheet
$win = $excel.ActiveWindow
$win.ScrollRow = 354
$win.ScrollColumn = 1
$ws.Range("D377").Select()
